$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$protectAddrs = @("D4", "D5", "D6", "D11", "D13", "D17", "D19", "D20", "D21", "D22", "D23", "D25", "D28", "D30", "D31", "D32", "D33", "D35", "D36", "D38", "D44", "D45", "D46", "D49", "D50", "D51")
foreach ($addr in $protectAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '68.173.90'
$ws.Range("E2").Value = '  +0.56%  '

$ws.Range("D3").Value = '3.798.43'
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '601.24'
$ws.Range("E5").Value = '  +0.71%  '

$ws.Range("D6").Value = '165.82'
$ws.Range("E6").Value = '  -0.91%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  -0.72%  '

$ws.Range("E9").Value = '  -1.24%  '

$ws.Range("E10").Value = '  +0.47%  '

$ws.Range("D11").Value = '6.50'
$ws.Range("E11").Value = '  +3.10%  '

$ws.Range("E12").Value = '  -1.17%  '

$ws.Range("D13").Value = '35.84'
$ws.Range("E13").Value = '  -0.79%  '

$ws.Range("D14").Value = '4.428.99'
$ws.Range("E14").Value = '  -0.39%  '

$ws.Range("D15").Value = '3.791.04'
$ws.Range("E15").Value = '  -0.78%  '

$ws.Range("D16").Value = '68.107.25'
$ws.Range("E16").Value = '  +0.48%  '

$ws.Range("D17").Value = '18.46'
$ws.Range("E17").Value = '  -1.02%  '

$ws.Range("E18").Value = '  +2.11%  '

$ws.Range("D19").Value = '7.09'
$ws.Range("E19").Value = '  -0.50%  '

$ws.Range("D20").Value = '461.64'
$ws.Range("E20").Value = '  +0.00%  '

$ws.Range("D21").Value = '9.73'
$ws.Range("E21").Value = '  -1.92%  '

$ws.Range("D22").Value = '0.700'
$ws.Range("E22").Value = '  -0.40%  '

$ws.Range("D23").Value = '0.0000150'
$ws.Range("E23").Value = '  -3.26%  '

$ws.Range("E24").Value = '  -0.89%  '

$ws.Range("D25").Value = '12.09'
$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("E26").Value = '  +0.56%  '

$ws.Range("E27").Value = '  -0.16%  '

$ws.Range("D28").Value = '10.00'
$ws.Range("E28").Value = '  -0.12%  '

$ws.Range("D29").Value = '3.942.92'
$ws.Range("E29").Value = '  -0.26%  '

$ws.Range("D30").Value = '7.41'
$ws.Range("E30").Value = '  +2.21%  '

$ws.Range("D31").Value = '2.65'
$ws.Range("E31").Value = '  -5.13%  '

$ws.Range("D32").Value = '2.23'
$ws.Range("E32").Value = '  -1.72%  '

$ws.Range("D33").Value = '29.38'
$ws.Range("E33").Value = '  -1.19%  '

$ws.Range("E34").Value = '  -0.10%  '

$ws.Range("D35").Value = '9.02'
$ws.Range("E35").Value = '  -0.75%  '

$ws.Range("D36").Value = '0.100'
$ws.Range("E36").Value = '  +0.11%  '

$ws.Range("E37").Value = '  +0.64%  '

$ws.Range("D38").Value = '3.28'
$ws.Range("E38").Value = '  -3.24%  '

$ws.Range("E39").Value = '  +0.05%  '

$ws.Range("E40").Value = '  -1.03%  '

$ws.Range("E41").Value = '  +0.02%  '

$ws.Range("E43").Value = '  +0.70%  '

$ws.Range("B44").Value = 'Arweave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D44").Value = '43.47'
$ws.Range("E44").Value = '  -0.89%  '

$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '47.39'
$ws.Range("E45").Value = '  -1.55%  '

$ws.Range("D46").Value = '151.43'
$ws.Range("E46").Value = '  +0.82%  '

$ws.Range("E47").Value = '  +0.34%  '

$ws.Range("E48").Value = '  +2.94%  '

$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").Value = '1.36'
$ws.Range("E49").Value = '  +6.40%  '

$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").Value = '393.94'
$ws.Range("E50").Value = '  +0.44%  '

$ws.Range("D51").Value = '26.73'
$ws.Range("E51").Value = '  +0.44%  '

foreach ($addr in $protectAddrs) {
    $ws.Range($addr).ClearFormats()
}
